$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 and B3 with their new numeric values (previously "NAN" text)
$ws.Range("B2").Value = 100
$ws.Range("B3").Value = 73.361080202582

# Update formulas in B4 and B5 to use 100 instead of 1523
$ws.Range("B4").Formula = "=B2/(12*100)"
$ws.Range("B5").Formula = "=B3/(12*100)"

# Update B6 and B7 values
$ws.Range("B6").Value = 0.086562964848485
$ws.Range("B7").Value = 16.918977604633
